$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "Outdated"
$ws.Range("D7").Value = "Outdated"
$ws.Range("D10").Value = "Outdated"

$ws.Range("D8").Select()
